# "fixed typo on question #1"
#
# Question #1 of the midterm gives an example array and a set of
# "N -> result" answers for findNthSmallest().  The array accidentally
# contained an extra "4, " element, and the worked example results were
# left over from an (incorrect) "Nth largest" computation instead of
# the "Nth smallest" the function actually implements.  This fixes
# both the array literal and the corresponding answer key line.

$d = $word.ActiveDocument

# Remove the stray "4, " from the example array: {1, 3, 2, 5, 4, 9, 8, 6} -> {1, 3, 2, 5, 9, 8, 6}
$d.Content.Find.Execute(
    "Given the array {1, 3, 2, 5, 4, 9, 8, 6}, the following",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Given the array {1, 3, 2, 5, 9, 8, 6}, the following",
    2)

# Correct the worked-example answer key to match findNthSmallest on the
# corrected array: N = 1 -> 9; N = 2 -> 9; N = 3 -> 6; N = 4 -> 5...
# becomes                   N = 1 -> 1; N = 2 -> 2; N = 3 -> 3; N = 4 -> 5...
$d.Content.Find.Execute(
    "N = 1 -> 9; N = 2 -> 9; N = 3 -> 6; N = 4 -> 5",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "N = 1 -> 1; N = 2 -> 2; N = 3 -> 3; N = 4 -> 5",
    2)
